$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Swap the M29:P29 <-> M30:P30 cell contents -------------------------
# (row 29 had the numeric alpha values, row 30 had the SES/Croston/SBA/SBJ
#  labels; the edit swaps them so the labels are on row 29 and the values
#  moved down to row 30)
$m29 = $ws.Range("M29").Value2
$n29 = $ws.Range("N29").Value2
$o29 = $ws.Range("O29").Value2
$p29 = $ws.Range("P29").Value2

$m30 = $ws.Range("M30").Value2
$n30 = $ws.Range("N30").Value2
$o30 = $ws.Range("O30").Value2
$p30 = $ws.Range("P30").Value2

$ws.Range("M29").Value2 = $m30
$ws.Range("N29").Value2 = $n30
$ws.Range("O29").Value2 = $o30
$ws.Range("P29").Value2 = $p30

$ws.Range("M30").Value2 = $m29
$ws.Range("N30").Value2 = $n29
$ws.Range("O30").Value2 = $o29
$ws.Range("P30").Value2 = $p29

# --- Chart title ----------------------------------------------------------
$chart = $ws.ChartObjects(1).Chart
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Croston's and Exponential Smoothing forecasts of item BIP001271, optimal alpha"

# --- Selection / view -------------------------------------------------
$ws.Range("E9").Select()
